# Daily attendance processing - 2025-11-24 16:54:09
#
# Column G ("Recorded By") holds a comma-separated list of recorder names.
# Normalize the ordering so that the literal token "System" (exact case)
# always comes first in the list, preserving the relative order of the
# remaining tokens. Rows whose "System" token is already first (or that
# contain no "System" token at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    if ($val -notlike "*System*") {
        continue
    }

    $parts = $val -split ", "

    $others = New-Object System.Collections.ArrayList
    $systemCount = 0
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $systemCount++
        } else {
            [void]$others.Add($p)
        }
    }

    if ($systemCount -eq 0) {
        continue
    }

    # Already ordered with System first and only one System token -> no change needed
    if (($parts[0].Equals("System")) -and ($systemCount -eq 1)) {
        continue
    }

    $newParts = New-Object System.Collections.ArrayList
    for ($i = 0; $i -lt $systemCount; $i++) {
        [void]$newParts.Add("System")
    }
    foreach ($p in $others) {
        [void]$newParts.Add($p)
    }

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
